$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wrapper element for MODS: the datastream-open tag (previously stored in Z1's
# matching opening fragment in C1) changes from an <update type="MODS"> wrapper to a
# <datastream type="md_descriptive" operation="update"> wrapper, and its closing
# counterpart in Z1 changes from </update> to </datastream> accordingly.
# Update Z1 first so the shared-string table ends up in the same append order as the
# authored workbook.
$ws.Range("Z1").Value = "</mods:mods></datastream></object>"
$ws.Range("C1").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink">'

# Restore the view to show column A with C1 selected (matches the saved view state).
$ws.Range("A1").Select()
$ws.Range("C1").Select()
